{"js": "// 1) \"nuestro banco de datos\" -> \"nuestra base de datos\" in the main\n//    certification paragraph (gender agreement fix: nuestro->nuestra,\n//    banco->base).\nconst oldPhrase = context.document.body.search(\"nuestro banco de datos\", { matchCase: true, matchWholeWord: false });\noldPhrase.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < oldPhrase.items.length; i++) {\n  oldPhrase.items[i].insertText(\"nuestra base de datos\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Update the cached result of the \"TIME \\@ \"d 'de' MMMM 'de' yyyy\"\"\n//    field (the date shown at the bottom of the certificate) from\n//    \"1 de noviembre de 2024\" to \"24 de febrero de 2025\".\nconst oldDate = context.document.body.search(\"1 de noviembre de 2024\", { matchCase: true, matchWholeWord: false });\noldDate.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < oldDate.items.length; i++) {\n  oldDate.items[i].insertText(\"24 de febrero de 2025\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"nuestro banco de datos\" -> \"nuestra base de datos\" in the main\n#    certification paragraph (gender agreement fix: nuestro->nuestra,\n#    banco->base).\n$find1 = $d.Content.Find\n$find1.Text = \"nuestro banco de datos\"\n$find1.MatchCase = $true\n$find1.Replacement.Text = \"nuestra base de datos\"\n$find1.Execute($null, $true, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# 2) Update the cached result of the \"TIME \\@ \"d 'de' MMMM 'de' yyyy\"\"\n#    field (the date shown at the bottom of the certificate) from\n#    \"1 de noviembre de 2024\" to \"24 de febrero de 2025\".\n$find2 = $d.Content.Find\n$find2.Text = \"1 de noviembre de 2024\"\n$find2.MatchCase = $true\n$find2.Replacement.Text = \"24 de febrero de 2025\"\n$find2.Execute($null, $true, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n"}
